# Refresh the cryptos list with the latest scraped price/volume snapshot.
# Numeric-looking "Price" strings are forced to text (NumberFormat "@")
# before assignment so Excel doesn't silently coerce/round them (e.g.
# "213.39" -> 213.39 as a float, losing the text formatting) the way the
# source data needs them stored (as plain text, matching the original
# inline-string cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.706.90"
$ws.Range("E2").Value = "  +1.32%  "
$ws.Range("D3").Value = "1.633.51"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.39"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.499"
$ws.Range("E6").Value = "  +3.02%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +1.36%  "
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.25"
$ws.Range("E10").Value = "  +2.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0841"
$ws.Range("E11").Value = "  +3.31%  "
$ws.Range("D12").Value = "1.859.41"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "1.642.14"
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.525"
$ws.Range("D16").Value = "26.671.96"
$ws.Range("E16").Value = "  +1.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.56"
$ws.Range("E17").Value = "  +1.64%  "
$ws.Range("E18").Value = "  +2.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.28"
$ws.Range("E19").Value = "  +7.90%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("E21").Value = "  +0.89%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.36"
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.16"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("E24").Value = "  +4.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.71"
$ws.Range("E25").Value = "  +2.20%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  +1.42%  "
$ws.Range("E28").Value = "  +3.88%  "
$ws.Range("E29").Value = "  +2.16%  "
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("E32").Value = "  +3.57%  "
$ws.Range("E33").Value = "  +2.18%  "
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.40"
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D36").Value = "1.223.76"
$ws.Range("E36").Value = "  +5.21%  "
$ws.Range("E37").Value = "  +5.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.808"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.501"
$ws.Range("E41").Value = "  -1.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.796"
$ws.Range("E42").Value = "  +1.63%  "
$ws.Range("E43").Value = "  -1.16%  "
$ws.Range("D44").Value = "1.767.12"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.72"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.57"
$ws.Range("E46").Value = "  +2.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.36"
$ws.Range("E47").Value = "  +2.80%  "
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("E49").Value = "  +0.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.65"
$ws.Range("E50").Value = "  +4.30%  "
$ws.Range("E51").Value = "  -0.26%  "
